$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B to hold the "Image" field,
# shifting Answer/Option1-4 columns one to the right.
$ws.Columns("B").Insert()

# New header + value for the inserted "Image" column.
$ws.Range("B1").Value = "Image"
$ws.Range("B2").Value = "img.png"

# Match the (non-bestFit) custom width recorded for the new column
# (closest value the character-width column model can reach to 23.140625).
$ws.Columns("B").ColumnWidth = 22.29

# Duplicate the single question row (row 2) five more times (rows 3-7),
# turning this into six identical question rows total.
for ($i = 3; $i -le 7; $i++) {
    $ws.Range("A2:G2").Copy()
    $ws.Range("A$i").PasteSpecial()
}

# Leave the same active cell/selection state recorded in the workbook.
$ws.Range("F14").Select()
